$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Extend data: rows 2-13 carry the refreshed NATMI LR-pair stats
# (Vcan -> Selp) across all 4x3 sending/target cluster combinations.

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Vcan"
$ws.Range("C2").Value = "Selp"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 3.204747333333333
$ws.Range("H2").Value = 9.614241999999999
$ws.Range("I2").Value = 0.01973032100547387
$ws.Range("J2").Value = 0.01973032100547387
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 15.960008
$ws.Range("N2").Value = 47.880024
$ws.Range("O2").Value = 0.9899803616776065
$ws.Range("P2").Value = 0.9899803616776066
$ws.Range("Q2").Value = 51.14779307797866
$ws.Range("R2").Value = 460.3301377018079
$ws.Range("S2").Value = 0.0195326303250143
$ws.Range("T2").Value = 0.01953263032501431

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Vcan"
$ws.Range("C3").Value = "Selp"
$ws.Range("D3").Value = "M2"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 3.204747333333333
$ws.Range("H3").Value = 9.614241999999999
$ws.Range("I3").Value = 0.01973032100547387
$ws.Range("J3").Value = 0.01973032100547387
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.06197766666666666
$ws.Range("N3").Value = 0.185933
$ws.Range("O3").Value = 0.003844401134548353
$ws.Range("P3").Value = 0.003844401134548354
$ws.Range("Q3").Value = 0.1986227619762222
$ws.Range("R3").Value = 1.787604857786
$ws.Range("S3").Value = 0.00007585126845844696
$ws.Range("T3").Value = 0.00007585126845844698

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Vcan"
$ws.Range("C4").Value = "Selp"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 3.204747333333333
$ws.Range("H4").Value = 9.614241999999999
$ws.Range("I4").Value = 0.01973032100547387
$ws.Range("J4").Value = 0.01973032100547387
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.09955433333333334
$ws.Range("N4").Value = 0.298663
$ws.Range("O4").Value = 0.006175237187845165
$ws.Range("P4").Value = 0.006175237187845166
$ws.Range("Q4").Value = 0.3190464842717778
$ws.Range("R4").Value = 2.871418358446
$ws.Range("S4").Value = 0.0001218394120011249
$ws.Range("T4").Value = 0.0001218394120011249

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Vcan"
$ws.Range("C5").Value = "Selp"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 145.2141163333334
$ws.Range("H5").Value = 435.6423490000001
$ws.Range("I5").Value = 0.8940240311559332
$ws.Range("J5").Value = 0.8940240311559333
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 15.960008
$ws.Range("N5").Value = 47.880024
$ws.Range("O5").Value = 0.9899803616776065
$ws.Range("P5").Value = 0.9899803616776066
$ws.Range("Q5").Value = 2317.618458392931
$ws.Range("R5").Value = 20858.56612553638
$ws.Range("S5").Value = 0.8850662337122225
$ws.Range("T5").Value = 0.8850662337122227

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Vcan"
$ws.Range("C6").Value = "Selp"
$ws.Range("D6").Value = "M2"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 145.2141163333334
$ws.Range("H6").Value = 435.6423490000001
$ws.Range("I6").Value = 0.8940240311559332
$ws.Range("J6").Value = 0.8940240311559333
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.06197766666666666
$ws.Range("N6").Value = 0.185933
$ws.Range("O6").Value = 0.003844401134548353
$ws.Range("P6").Value = 0.003844401134548354
$ws.Range("Q6").Value = 9.000032097401888
$ws.Range("R6").Value = 81.000288876617
$ws.Range("S6").Value = 0.003436986999689362
$ws.Range("T6").Value = 0.003436986999689363

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Vcan"
$ws.Range("C7").Value = "Selp"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 145.2141163333334
$ws.Range("H7").Value = 435.6423490000001
$ws.Range("I7").Value = 0.8940240311559332
$ws.Range("J7").Value = 0.8940240311559333
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.3333333333333333
$ws.Range("M7").Value = 0.09955433333333334
$ws.Range("N7").Value = 0.298663
$ws.Range("O7").Value = 0.006175237187845165
$ws.Range("P7").Value = 0.006175237187845166
$ws.Range("Q7").Value = 14.45669454215411
$ws.Range("R7").Value = 130.110250879387
$ws.Range("S7").Value = 0.005520810444021363
$ws.Range("T7").Value = 0.005520810444021364

# Row 8
$ws.Range("A8").Value = "M2"
$ws.Range("B8").Value = "Vcan"
$ws.Range("C8").Value = "Selp"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.8052786666666667
$ws.Range("H8").Value = 2.415836
$ws.Range("I8").Value = 0.004957771998726471
$ws.Range("J8").Value = 0.004957771998726472
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 15.960008
$ws.Range("N8").Value = 47.880024
$ws.Range("O8").Value = 0.9899803616776065
$ws.Range("P8").Value = 0.9899803616776066
$ws.Range("Q8").Value = 12.85225396222933
$ws.Range("R8").Value = 115.670285660064
$ws.Range("S8").Value = 0.004908096916414342
$ws.Range("T8").Value = 0.004908096916414344

# Row 9
$ws.Range("A9").Value = "M2"
$ws.Range("B9").Value = "Vcan"
$ws.Range("C9").Value = "Selp"
$ws.Range("D9").Value = "M2"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.8052786666666667
$ws.Range("H9").Value = 2.415836
$ws.Range("I9").Value = 0.004957771998726471
$ws.Range("J9").Value = 0.004957771998726472
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.06197766666666666
$ws.Range("N9").Value = 0.185933
$ws.Range("O9").Value = 0.003844401134548353
$ws.Range("P9").Value = 0.003844401134548354
$ws.Range("Q9").Value = 0.04990929277644444
$ws.Range("R9").Value = 0.449183634988
$ws.Range("S9").Value = 0.0000190596642967361
$ws.Range("T9").Value = 0.00001905966429673611

# Row 10
$ws.Range("A10").Value = "M2"
$ws.Range("B10").Value = "Vcan"
$ws.Range("C10").Value = "Selp"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.8052786666666667
$ws.Range("H10").Value = 2.415836
$ws.Range("I10").Value = 0.004957771998726471
$ws.Range("J10").Value = 0.004957771998726472
$ws.Range("K10").Value = 1
$ws.Range("L10").Value = 0.3333333333333333
$ws.Range("M10").Value = 0.09955433333333334
$ws.Range("N10").Value = 0.298663
$ws.Range("O10").Value = 0.006175237187845165
$ws.Range("P10").Value = 0.006175237187845166
$ws.Range("Q10").Value = 0.08016898080755556
$ws.Range("R10").Value = 0.721520827268
$ws.Range("S10").Value = 0.00003061541801539316
$ws.Range("T10").Value = 0.00003061541801539316

# Row 11
$ws.Range("A11").Value = "sCs"
$ws.Range("B11").Value = "Vcan"
$ws.Range("C11").Value = "Selp"
$ws.Range("D11").Value = "ECs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 13.203389
$ws.Range("H11").Value = 39.610167
$ws.Range("I11").Value = 0.08128787583986632
$ws.Range("J11").Value = 0.08128787583986634
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 15.960008
$ws.Range("N11").Value = 47.880024
$ws.Range("O11").Value = 0.9899803616776065
$ws.Range("P11").Value = 0.9899803616776066
$ws.Range("Q11").Value = 210.726194067112
$ws.Range("R11").Value = 1896.535746604008
$ws.Range("S11").Value = 0.08047340072395523
$ws.Range("T11").Value = 0.08047340072395526

# Row 12
$ws.Range("A12").Value = "sCs"
$ws.Range("B12").Value = "Vcan"
$ws.Range("C12").Value = "Selp"
$ws.Range("D12").Value = "M2"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 13.203389
$ws.Range("H12").Value = 39.610167
$ws.Range("I12").Value = 0.08128787583986632
$ws.Range("J12").Value = 0.08128787583986634
$ws.Range("K12").Value = 2
$ws.Range("L12").Value = 0.6666666666666666
$ws.Range("M12").Value = 0.06197766666666666
$ws.Range("N12").Value = 0.185933
$ws.Range("O12").Value = 0.003844401134548353
$ws.Range("P12").Value = 0.003844401134548354
$ws.Range("Q12").Value = 0.8183152423123333
$ws.Range("R12").Value = 7.364837180811
$ws.Range("S12").Value = 0.0003125032021038078
$ws.Range("T12").Value = 0.0003125032021038079

# Row 13
$ws.Range("A13").Value = "sCs"
$ws.Range("B13").Value = "Vcan"
$ws.Range("C13").Value = "Selp"
$ws.Range("D13").Value = "sCs"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 13.203389
$ws.Range("H13").Value = 39.610167
$ws.Range("I13").Value = 0.08128787583986632
$ws.Range("J13").Value = 0.08128787583986634
$ws.Range("K13").Value = 1
$ws.Range("L13").Value = 0.3333333333333333
$ws.Range("M13").Value = 0.09955433333333334
$ws.Range("N13").Value = 0.298663
$ws.Range("O13").Value = 0.006175237187845165
$ws.Range("P13").Value = 0.006175237187845166
$ws.Range("Q13").Value = 1.314454589635667
$ws.Range("R13").Value = 11.830091306721
$ws.Range("S13").Value = 0.000501971913807283
$ws.Range("T13").Value = 0.0005019719138072832

